$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.334.74"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3
$ws.Range("D3").Value = "'1.846.11"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4
$ws.Range("D4").Value = "'0.9976"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'239.98"
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("E6").Value = "  -0.61%  "

# Row 7
$ws.Range("D7").Value = "'0.9984"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").Value = "'0.07596"
$ws.Range("E8").Value = "  -1.30%  "

# Row 9
$ws.Range("D9").Value = "'0.2899"
$ws.Range("E9").Value = "  -1.49%  "

# Row 10
$ws.Range("E10").Value = "  +0.83%  "

# Row 11
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("D12").Value = "'5.022"
$ws.Range("E12").Value = "  -0.13%  "

# Row 13
$ws.Range("D13").Value = "'0.6779"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14
$ws.Range("E14").Value = "  -2.29%  "

# Row 15
$ws.Range("D15").Value = "'82.90"

# Row 16
$ws.Range("D16").Value = "'6.125"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17
$ws.Range("D17").Value = "'29.372.18"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("D18").Value = "'227.55"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("E19").Value = "  -1.14%  "

# Row 20
$ws.Range("D20").Value = "'0.9983"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("D21").Value = "'7.462"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$ws.Range("D22").Value = "'0.9982"
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("D23").Value = "'158.53"
$ws.Range("E23").Value = "  +0.79%  "

# Row 24
$ws.Range("D24").Value = "'0.1381"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
$ws.Range("D25").Value = "'8.423"
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
$ws.Range("D26").Value = "'17.63"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("D27").Value = "'1.421"
$ws.Range("E27").Value = "  +8.19%  "

# Row 28
$ws.Range("D28").Value = "'1.458"
$ws.Range("E28").Value = "  -0.72%  "

# Row 29
$ws.Range("D29").Value = "'0.05595"
$ws.Range("E29").Value = "  -2.55%  "

# Row 30
$ws.Range("E30").Value = "  -0.35%  "

# Row 31
$ws.Range("D31").Value = "'4.065"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").Value = "'1.830"
$ws.Range("E33").Value = "  -1.06%  "

# Row 34
$ws.Range("D34").Value = "'0.6919"
$ws.Range("E34").Value = "  -2.43%  "

# Row 35
$ws.Range("D35").Value = "'2.581"
$ws.Range("E35").Value = "  -0.28%  "

# Row 36
$ws.Range("D36").Value = "'0.01798"

# Row 37
$ws.Range("D37").Value = "'1.225.37"

# Row 38
$ws.Range("D38").Value = "'2.723"
$ws.Range("E38").Value = "  -1.91%  "

# Row 39
$ws.Range("D39").Value = "'6.353"
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").Value = "'0.8958"

# Row 41
$ws.Range("D41").Value = "'0.9981"

# Row 42
$ws.Range("D42").Value = "'101.28"
$ws.Range("E42").Value = "  -0.55%  "

# Row 43
$ws.Range("D43").Value = "'65.44"
$ws.Range("E43").Value = "  -1.32%  "

# Row 44
$ws.Range("D44").Value = "'7.190"
$ws.Range("E44").Value = "  +0.50%  "

# Row 45
$ws.Range("D45").Value = "'0.3985"

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'9.007"
$ws.Range("E46").Value = "  -0.35%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.688"
$ws.Range("E47").Value = "  -0.04%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1140"
$ws.Range("E48").Value = "  +1.32%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000113"
$ws.Range("E49").Value = "  -6.23%  "

# Row 50
$ws.Range("D50").Value = "'0.05693"
$ws.Range("E50").Value = "  -0.41%  "

# Row 51
$ws.Range("E51").Value = "  -0.19%  "
